$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "dSF" (column F) values after re-pulling data.
$updates = @{
    2  = -5
    11 = 1
    15 = 0
    18 = 0
    20 = 1
    21 = 0
    29 = 3
    33 = -8
    35 = 1
    37 = 0
    39 = 1
    46 = 1
    47 = 2
    48 = 0
    51 = -1
    52 = 1
    53 = 0
    56 = 1
    58 = 1
    68 = 1
    72 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
